# Insert a new data row right after the header block of existing rows,
# at worksheet row 439 (pushing the former rows 439..541 down to 440..542).
# This mirrors the OOXML diff: dimension grows from A1:R541 to A1:R542 and
# a brand-new "Paine / 1a (guarda)" record is inserted in the middle of the
# "Zapallo" table, with every row after it shifted down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 439 downward by inserting a blank row at position 439.
$ws.Rows("439").Insert()

# Populate the newly inserted row 439 with the new record's values.
$ws.Cells.Item(439, 1).Value  = 4
$ws.Cells.Item(439, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(439, 3).Value  = "Los Lagos"
$ws.Cells.Item(439, 4).Value  = 45204
$ws.Cells.Item(439, 5).Value  = 10
$ws.Cells.Item(439, 6).Value  = 100112045
$ws.Cells.Item(439, 7).Value  = "Zapallo"
$ws.Cells.Item(439, 8).Value  = "Paine"
$ws.Cells.Item(439, 9).Value  = "1a (guarda)"
$ws.Cells.Item(439, 10).Value = 750
$ws.Cells.Item(439, 11).Value = 800
$ws.Cells.Item(439, 12).Value = 800
$ws.Cells.Item(439, 13).Value = 800
$ws.Cells.Item(439, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(439, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(439, 16).Value = 800
$ws.Cells.Item(439, 17).Value = 1
$ws.Cells.Item(439, 18).Value = "Hortaliza"
